$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Engine" task list is being edited:
#   - "Complete Renderer refactoring" (row 2) is removed entirely.
#   - "Deprecate Vector3. Replace it with position, direction and unit
#     direction" (originally row 4) is removed entirely - Vector3 is gone.
#   - A new task "Move to pre-compiled shaders" is added in its place.
# Net effect: one row shorter overall, and the two cell comments that lived
# lower in the sheet shift up by one row.
# ---------------------------------------------------------------------------

# Grab the two existing comments' text before we start shuffling rows, then
# remove them - their anchor cells are about to move.
$commentOneText = $ws.Range("B13").Comment.Text()
$commentTwoText = $ws.Range("B16").Comment.Text()
$ws.Range("B13").Comment.Delete()
$ws.Range("B16").Comment.Delete()

# Remove the "Complete Renderer refactoring" row entirely.
$ws.Rows(2).Delete()

# After the delete above, "Deprecate Vector3..." is now row 3 - remove it too.
$ws.Rows(3).Delete()

# Insert a fresh row 3 for the replacement task.
$ws.Rows(3).Insert()
$ws.Range("A3").Value = "Engine"
$ws.Range("B3").Value = "Move to pre-compiled shaders"
$ws.Range("C3").Value = 5

# Re-create the two comments at their new (shifted-up-by-one) locations.
$newCommentOne = $ws.Range("B12").AddComment($commentOneText)
$newCommentTwo = $ws.Range("B15").AddComment($commentTwoText)

# Match the final cursor/selection position recorded in the workbook.
$ws.Range("B20").Select()
